$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to make edits, then restore protection after.
$ws.Unprotect()

# Update the confidential/disclosure note date (2021-07-07 -> 2021-07-08)
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2483720611986374
$ws.Range("E2").Value = -0.01357575757575757

$ws.Range("D3").Value = 0.250143163153187
$ws.Range("E3").Value = -0.02004943696786576

$ws.Range("D4").Value = 0.2580014058431259
$ws.Range("E4").Value = -0.01393728222996515

$ws.Range("D5").Value = 0.2434833698050496
$ws.Range("E5").Value = -0.005763688760806795

$ws.Range("D6").Value = 0.9999999999999999
$ws.Range("E6").Value = -0.01338626924495856

# Re-protect the sheet to restore original protection state.
$ws.Protect()
